$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("46cea3b9-bf3a-400c-a9ca-f07bc02c51f0", "Waste", "paper", 10, 0, "2024-09-21", "17:23:07"),
    @("c0d653bf-2f1a-4cfa-988f-9ee6b267ef45", "Waste", "paper", 10, 0, "2024-09-21", "17:45:11"),
    @("229acf0c-779f-434d-a864-34fde3ed54e9", "Out", "paper", 100, 0, "2024-09-21", "17:45:52"),
    @("0535ef63-cd48-482b-806e-b8e7360c32ac", "Out", "paper", 1, 0, "2024-09-21", "17:51:43"),
    @("1899b5d8-72a6-4f60-91fc-a9c68866dd90", "Out", "paper", 2, 0, "2024-09-21", "17:52:14"),
    @("006ffc8f-2e6b-4a91-b27d-4b7a650f632d", "Out", "paper", 1234, 0, "2024-09-21", "17:54:32"),
    @("3d7a2e09-c12c-49df-86b9-f90c2a8fd76f", "Out", "Ram", 10, 10101, "2024-09-21", "17:56:11")
)

$row = 3
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]

    # Columns F (Date) and G (Time) hold plain text like "2024-09-21" /
    # "17:23:07". Assigning that text straight to .Value while the cell is
    # still General-formatted makes Excel auto-recognize it as a real
    # date/time and silently rewrite the stored value as a serial number.
    # Entering it as a formula that evaluates to the exact same text, then
    # collapsing the formula to its literal value via copy / paste-values,
    # keeps the literal string without leaving any extra number formatting
    # applied to the cell.
    $ws.Cells.Item($row, 6).Formula = "=""" + $entry[5] + """"
    $ws.Cells.Item($row, 7).Formula = "=""" + $entry[6] + """"
    $ws.Range($ws.Cells.Item($row, 6), $ws.Cells.Item($row, 7)).Copy()
    $ws.Range($ws.Cells.Item($row, 6), $ws.Cells.Item($row, 7)).PasteSpecial(-4163)

    $row++
}
